$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267, shifting existing rows 267:332 down to 268:333
# (dimension grows from A1:R332 to A1:R333).
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with a new price record.
$ws.Range("A267").Value = 5
$ws.Range("B267").Value = "Macroferia Regional de Talca"
$ws.Range("C267").Value = "Maule"
$ws.Range("D267").Value = 44722
$ws.Range("E267").Value = 7
$ws.Range("F267").Value = 100114014
$ws.Range("G267").Value = "Betarraga"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 3000
$ws.Range("K267").Value = 700
$ws.Range("L267").Value = 700
$ws.Range("M267").Value = 700
$ws.Range("N267").Value = "`$/paquete 5 unidades"
$ws.Range("O267").Value = "Región del Maule"
$ws.Range("P267").Value = 140
$ws.Range("Q267").Value = 5
$ws.Range("R267").Value = "Hortaliza"
